$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-74 down to 6-75.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new record. The row copies the
# same static columns (A,B,C,E,F,G,H,I,Q,R) as the record that used to sit
# at row 5 (now at row 6), but carries its own date / volume / price data.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C5").Value = "Los Lagos"
$ws.Range("D5").Value = 44691
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 100112031
$ws.Range("G5").Value = "Poroto verde"
$ws.Range("H5").Value = "Magnum"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 30000
$ws.Range("L5").Value = 30000
$ws.Range("M5").Value = 30000
$ws.Range("N5").Value = '$/saco 25 kilos'
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 1200
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"
